$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (this shifts B..J to D..L,
# and Excel automatically rewrites formula references accordingly)
$ws.Range("B1:C1").EntireColumn.Insert()

# ---- Row 1 headers (order matters for shared-string table layout) ----
$ws.Range("D1").Value = "bid_price"
$ws.Range("B1").Value = "atb_capex"
$ws.Range("C1").Value = "atb_capex_multiplier"

# ---- New atb_capex / atb_capex_multiplier values (column B / C) for rows 2-7 ----
$ws.Range("B2").Value = 1483
$ws.Range("C2").Value = 0.3

$ws.Range("B3").Value = 2478
$ws.Range("C3").Value = 0.3

$ws.Range("B4").Value = 2911
$ws.Range("C4").Value = 0.3

$ws.Range("B5").Value = 2000.5
$ws.Range("C5").Value = 0.3

$ws.Range("B6").Value = 1556
$ws.Range("C6").Value = 0.3

$ws.Range("B7").Value = 6318.5
$ws.Range("C7").Value = 0.3

# ---- fuel_type (column A) relabeling ----
$ws.Range("A6").Value = "gas"
$ws.Range("A7").Value = "wind_offshore"
$ws.Range("A8").Value = ""

# ---- sheet view tweaks (best-effort; scroll position to column G, select H9) ----
$wb.Windows.Item(1).ScrollColumn = 7
$excel.Goto($ws.Range("H9"), $false)
$ws.Range("H9").Select()
